$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

# New naive QoQ error table data (values shifted down one row, with a
# freshly computed row inserted at the top of the data block, row 2).
$data = @(
    @(-0.02907897629796788, 0.3131278957257717,  0.181524606355785,   0.4260570458938391, 0.43998257208981,   15),
    @( 0.0370360702193368,  0.2877205798324301,  0.135345957066826,   0.3678939481247633, 0.3798420736582994, 14),
    @( 0.01989109342689383, 0.3298745939743749,  0.1503185731831121,  0.3877093926939507, 0.4030092977225764, 13),
    @( 0.06273519545459039, 0.2460334154802962,  0.09558537405740077, 0.3091688439306276, 0.3161984462764572, 12),
    @( 0.02711121047383215, 0.2779584568479502,  0.1216586105046287,  0.34879594393374,   0.3647135205176791, 11),
    @(-0.01088135635153479, 0.2849421119723689,  0.09825224066429059, 0.3134521345664926, 0.3302084135617004, 10),
    @(-0.04428645741563344, 0.3765628100937468,  0.1690784592812682,  0.4111915116843588, 0.4335975266214011, 9),
    @(-0.01393931246739222, 0.3509192590318558,  0.1607817096315818,  0.4009759464501354, 0.4284022211487281, 8),
    @(-0.03688841855209302, 0.2799590153621541,  0.1329006694222477,  0.3645554408073589, 0.3917438359423185, 7),
    @(-0.06364682135181432, 0.1670219060428917,  0.03807752624502202, 0.1951346362002964, 0.2020691219662191, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

$wb.Save()
